$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the Heading1
#    title paragraph. We build it by cloning the trailing duplicate
#    bold-title paragraph (same empty-run + bold-run shape we need),
#    pasting it into place, then editing its runs' text in place so the
#    run/formatting structure (leading empty run, bold run, plain run)
#    matches exactly.
# ---------------------------------------------------------------------

$titleDup = $d.Paragraphs(49)      # "Play Bronco Spirit Free: Detailed Review" (bold, near the end)
$titleDup.Range.Copy()

$heading1 = $d.Paragraphs(1)       # "Play Bronco Spirit Free: Detailed Review" (Heading1)
$heading1.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs(2)
$newPara.Style = "Normal"          # strip the inherited Heading1 style -> plain paragraph
$newPara.Range.Paste()

$metaPara = $d.Paragraphs(2)
$metaRange = $metaPara.Range

# The pasted run currently reads "Play Bronco Spirit Free: Detailed Review" (bold).
# Replace just that run's text with "Meta description", keeping its bold formatting.
$oldLabelLen = [string]"Play Bronco Spirit Free: Detailed Review".Length
$labelRange = $d.Range($metaRange.Start, $metaRange.Start + $oldLabelLen)
$labelRange.Text = "Meta description"

# Append the rest of the sentence as a new, non-bold run right after the label.
$metaPara2 = $d.Paragraphs(2)
$metaRange2 = $metaPara2.Range
$insertPoint = $d.Range($metaRange2.End - 1, $metaRange2.End - 1)
$insertPoint.InsertAfter(": Find out why you should play Bronco Spirit, a Native American themed slot with high RTP and well-crafted graphics. Play now for free.")

# ---------------------------------------------------------------------
# 2) Remove the old duplicated bold title paragraph near the end of the
#    document (it shifted down by one slot because of the insertion
#    above).
# ---------------------------------------------------------------------

$dupTitle = $d.Paragraphs(50)
$dupTitle.Range.Delete()

# ---------------------------------------------------------------------
# 3) Replace the italic meta-description paragraph's text with the
#    DALLE image prompt, keeping its italic formatting intact.
# ---------------------------------------------------------------------

$descPara = $d.Paragraphs(50)
$descRange = $descPara.Range
$descTextRange = $d.Range($descRange.Start, $descRange.End - 1)
$descTextRange.Text = "Prompt for DALLE: Create a feature image for Bronco Spirit that showcases the game's Native American theme. The image should be in a cartoon style and prominently feature a happy Maya warrior with glasses. The warrior should be depicted riding a majestic horse against a desert background, with barren mountains, cactus plants, and vultures in the foreground. The overall color scheme should reflect the game's dominant orange hue, connoting the setting sun. Make sure the image highlights the essential elements of the game, such as the reels, the mustangs (golden coins), the bonus symbol (the sunset), and the Wild symbol (the gold coin depicting the horse)."
